$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Id" column (A) held integers that were breaking downstream int-parsing
# logic, so remove it entirely. This shifts UniqueId/Name/Email from B:D to A:C.
$ws.Columns("A").Delete()

# Update the lingering selection to land on the new first column.
$ws.Range("A8").Select() | Out-Null
